$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: B25 was stored as text "3"; change it to a real number 3.
$ws.Cells.Item(25, 2).Value = 3

# Row 26: new annotation row appended after row 25.
$ws.Cells.Item(26, 1).Value = "Ruilin"
# B26 should remain a text value "3" (mirrors the previous row 25 layout).
$ws.Cells.Item(26, 2).Value = "'3"
$ws.Cells.Item(26, 2).Style = "Normal"
$ws.Cells.Item(26, 3).Value = "无"
$ws.Cells.Item(26, 4).Value = "QSN"
$ws.Cells.Item(26, 5).Value = "MET"
$ws.Cells.Item(26, 6).Value = "f1a2d8e0-a083-4e7a-9e83-7f61c3c0d7bb"
$ws.Cells.Item(26, 7).Value = "HJewuJWCZ_annotated.xlsx"
$ws.Cells.Item(26, 8).Value = "Is accuracy stable, can it drop back down below the threshold in the next epoch?"
